$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (A1:F1) with the new column order
$ws.Cells.Item(1, 1).Value = "kitchens_1"
$ws.Cells.Item(1, 2).Value = "bedrooms_1"
$ws.Cells.Item(1, 3).Value = "bedrooms_2"
$ws.Cells.Item(1, 4).Value = "kitchens_2"
$ws.Cells.Item(1, 5).Value = "living_rooms_1"
$ws.Cells.Item(1, 6).Value = "living_rooms_2"

# Update the 0/1 matrix rows (A2:F7) to match the new column assignment
$matrix = @(
    @(0,1,0,0,0,0),
    @(1,0,0,0,0,0),
    @(0,0,0,0,0,1),
    @(0,0,1,0,0,0),
    @(0,0,0,1,0,0),
    @(0,0,0,0,1,0)
)

for ($i = 0; $i -lt 6; $i++) {
    $rowVals = $matrix[$i]
    for ($j = 0; $j -lt 6; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $rowVals[$j]
    }
}
